$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 2 (gab.rioja@gab.es / RE: 202540501668) is removed entirely.
# This shifts the old row 3 (SAP España / SAP Business Suite Innovation Day)
# up into row 2, and the used range shrinks from A1:C3 to A1:C2.
$ws.Rows(2).Delete()
